$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='30.395.10'; E='  -1.50%  '},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.915.27'; E='  +1.69%  '},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.001'; E='  +0.32%  '},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='241.44'; E='  +1.60%  '},
    @{Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.001'; E='  +0.25%  '},
    @{Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4682'; E='  -2.31%  '},
    @{Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2837'; E='  -0.26%  '},
    @{Row=9; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.06967'; E='  +6.78%  '},
    @{Row=10; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='105.55'; E='  +10.63%  '},
    @{Row=11; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='18.08'; E='  -3.74%  '},
    @{Row=12; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.900.36'; E='  +0.93%  '},
    @{Row=13; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07639'; E='  +1.61%  '},
    @{Row=14; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.156'; E='  +0.88%  '},
    @{Row=15; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6535'; E='  -0.05%  '},
    @{Row=16; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='284.24'; E='  -4.35%  '},
    @{Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='30.382.79'; E='  -1.25%  '},
    @{Row=18; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.001'; E='  +0.04%  '},
    @{Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007598'; E='  +0.48%  '},
    @{Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.95'; E='  -1.49%  '},
    @{Row=21; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.147.09'; E='  +1.44%  '},
    @{Row=22; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.001'; E='  +0.05%  '},
    @{Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.200'; E='  +0.59%  '},
    @{Row=24; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='6.158'; E='  +0.56%  '},
    @{Row=25; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='167.96'; E='  -0.61%  '},
    @{Row=26; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='9.201'; E='  -0.64%  '},
    @{Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='20.86'; E='  +6.54%  '},
    @{Row=28; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.028'; E='  +3.07%  '},
    @{Row=29; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.1066'; E='  +1.11%  '},
    @{Row=30; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.368'; E='  -0.26%  '},
    @{Row=31; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.104'; E='  -0.95%  '},
    @{Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.947'; E='  -0.18%  '},
    @{Row=33; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05055'; E='  +1.24%  '},
    @{Row=34; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7354'; E='  +1.71%  '},
    @{Row=35; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.141'; E='  -3.16%  '},
    @{Row=36; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.000'; E='  +0.23%  '},
    @{Row=37; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.716'; E='  +0.56%  '},
    @{Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01996'; E='  +3.24%  '},
    @{Row=39; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.666'; E='  -1.90%  '},
    @{Row=40; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.044'; E='  -1.14%  '},
    @{Row=41; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='108.10'; E='  +0.50%  '},
    @{Row=42; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8702'; E='  -2.61%  '},
    @{Row=43; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.805'; E='  +3.71%  '},
    @{Row=44; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.001'; E='  +0.32%  '},
    @{Row=45; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4180'; E='  -0.81%  '},
    @{Row=46; B='BitcoinSV'; C='https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'; D='51.95'; E='  +23.33%  '},
    @{Row=47; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='66.98'; E='  +2.29%  '},
    @{Row=48; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.106'; E='  -3.38%  '},
    @{Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.124'; E='  +3.06%  '},
    @{Row=50; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1200'; E='  -2.61%  '},
    @{Row=51; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='34.40'; E='  -0.74%  '}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
